# Edit script: add X3/Y3 values, and append two new rows (4 and 5)
# reflecting additional "repeater" trade runs (ran trades overnight).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in PriceChange (X) and UpDown (Y) ---
$ws.Range("X3").Value = -0.34999899999999684
$ws.Range("Y3").Value = "Down"

# --- Row 4: new trade run (copy formatting from row 3, then set values) ---
$ws.Range("A3:Y3").Copy($ws.Range("A4:Y4"))
$ws.Range("A4").Value = 42649.612187500003
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 38
$ws.Range("E4").Value = 5843
$ws.Range("F4").Value = 294
$ws.Range("G4").Value = 62
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 91
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 9503
$ws.Range("L4").Value = 56
$ws.Range("M4").Value = 34
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = "Named"
$ws.Range("Q4").Value = 35.483823948801813
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.0965
$ws.Range("T4").Value = 0.0269
$ws.Range("U4").Value = 4.82
$ws.Range("V4").Value = 2.2799999999999998
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = -0.34999899999999684
$ws.Range("Y4").Value = "Down"

# --- Row 5: new trade run (copy formatting from row 4, then set values) ---
$ws.Range("A4:Y4").Copy($ws.Range("A5:Y5"))
$ws.Range("A5").Value = 42649.635555555556
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 56
$ws.Range("E5").Value = 7286
$ws.Range("F5").Value = 422
$ws.Range("G5").Value = 67
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 91
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 11425
$ws.Range("L5").Value = 85
$ws.Range("M5").Value = 41
$ws.Range("N5").Value = 44
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = "Named"
$ws.Range("Q5").Value = 35.483823948801813
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.0965
$ws.Range("T5").Value = 0.0269
$ws.Range("U5").Value = 4.82
$ws.Range("V5").Value = 2.2799999999999998
$ws.Range("W5").Value = 0
# Row 5 has no PriceChange/UpDown recorded yet (trade still open)
$ws.Range("X5").ClearContents()
$ws.Range("Y5").ClearContents()
